$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text string into a cell without Excel auto-converting
# a numeric-looking string (e.g. "0.0000") into a real number, and without
# disturbing the cell's existing number format / style (it should stay
# General, centered, like the other format-name header labels on row 1).
# We do this by computing the text with a formula on a scratch cell (a
# formula result is always stored with its native type, text in this case,
# regardless of the cell's number format), copying just that computed value
# onto the destination cell, then wiping the scratch cell again.
$staging = $ws.Cells.Item(60, 1)
function Set-TextValue($cell, [string]$text) {
    $staging.Formula = '="' + $text + '"'
    $staging.Copy()
    $cell.PasteSpecial(-4163) # xlPasteValues
}

# ------------------------------------------------------------------
# 1. Generalize the number-format columns D,E,F (previously thousands
#    separator formats) to plain decimal-precision formats.
# ------------------------------------------------------------------
Set-TextValue $ws.Cells.Item(1, 4) "0.0000"
Set-TextValue $ws.Cells.Item(1, 5) "0.00000"
Set-TextValue $ws.Cells.Item(1, 6) "0.0000#"
$staging.Clear()

$ws.Cells.Item(2, 4).NumberFormat = "0.0000"
$ws.Cells.Item(2, 5).NumberFormat = "0.00000"
$ws.Cells.Item(2, 6).NumberFormat = "0.0000#"

# ------------------------------------------------------------------
# 2. Add two new generalized date/time format demo columns, J and K.
# ------------------------------------------------------------------
$ws.Cells.Item(1, 10).Value = "mmddyy"
$ws.Cells.Item(1, 11).Value = "mmddyyam/pmdd"

$ws.Cells.Item(2, 10).Value = 51276
$ws.Cells.Item(2, 10).NumberFormat = "MMDDYY"

$ws.Cells.Item(2, 11).Value = 51276
$ws.Cells.Item(2, 11).NumberFormat = "MMDDYYAM/PMDD"

# ------------------------------------------------------------------
# 3. Widen the new K column to fit the longer format label.
#    (ColumnWidth uses Excel's character-width units and gets rounded to
#    the nearest whole pixel internally; 20 is the input that lands on
#    the target stored width of 20.83 for this column position.)
# ------------------------------------------------------------------
$ws.Columns.Item(11).ColumnWidth = 20

# ------------------------------------------------------------------
# 4. Move the active selection, matching the saved view state.
# ------------------------------------------------------------------
$ws.Range("I5").Select()
